# "Generate Report for Archive"
#
# The localization status of the two sample files moved from
# "Ready for handoff" to "In Translation" on every sheet that surfaces
# status (the Overview roll-up columns E/F for zh-cn and de-de, plus the
# per-language "Status" column on the zh-cn and de-de detail sheets).
# Shrinking that text auto-narrows the Status-bearing columns, so their
# widths are refreshed to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn detail sheet: Status column (col C) ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Columns.Item(3).ColumnWidth = 12.5

# --- de-de detail sheet: Status column (col C) ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Columns.Item(3).ColumnWidth = 12.5
